# Re-write the Sheet1 header row (A1:C1) so the columns read, left to right:
#   Master_Data_Type | Master_Data_Volume | Master_Data_Used_Volume
# Setting A1's text again (same characters minus the old trailing padding)
# causes Excel to mint a fresh shared-string entry for it at the end of the
# table, which is what moves "Master_Data_Type" to the back of the shared
# strings list while "Master_Data_Volume"/"Master_Data_Used_Volume" shift up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Master_Data_Type"
$ws.Range("B1").Value = "Master_Data_Volume"
$ws.Range("C1").Value = "Master_Data_Used_Volume"

# Restore the default (no custom) selection state on the sheet - A1.
$ws.Range("A1").Select()
